$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 11 (FUENTES PAREDES MARIA FERNANDA)
$ws1.Range("D11").Value = 2595.84
$ws1.Range("G11").Value = 43.54
$ws1.Range("L11").Value = 1265.09

# Row 14 (HERRERA CAICEDO LUIS FRANKLIN)
$ws1.Range("D14").Value = 2092.04

# Row 32 (summary counts "X de 30")
$ws1.Range("D32").Value = "4 de 30"
$ws1.Range("G32").Value = "1 de 30"
$ws1.Range("L32").Value = "3 de 30"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F11").Value = 3997
$ws2.Range("F14").Value = 3319.36
$ws2.Range("F32").Value = 13450.89

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# widen column E slightly (22 -> 23 in stored units)
$ws3.Columns.Item(5).ColumnWidth = 22.17

# Row 3 (240X80 PORCELANATO)
$ws3.Range("D3").Value = 5853.51
$ws3.Range("E3").Value = -2733.3955
$ws3.Range("F3").Value = 1.876056151144453

# Row 6 (GRIFERIAS)
$ws3.Range("D6").Value = 43.54
$ws3.Range("E6").Value = 63.27999999999999
$ws3.Range("F6").Value = 0.4076015727391875

# Row 15 (PIEDRA SINTERIZADA)
$ws3.Range("D15").Value = 2132.74
$ws3.Range("E15").Value = -1605.71
$ws3.Range("F15").Value = 4.046714608276568

# Row 19 (TOTAL)
$ws3.Range("D19").Value = 13445.13
$ws3.Range("E19").Value = 16092.66107555787
$ws3.Range("F19").Value = 0.4551840036246199
